# filter search by year
# Adds the missing "Queen of Clubs" (QC) card entries to the Languedoc
# card data sheet, mirroring the pattern used for the other existing
# suit/rank combinations (Recto + Verso rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row after the current data (right after row 23).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow1 = $lastRow + 1
$newRow2 = $lastRow + 2

# Row: Recto (1) - Queen of Clubs
$ws.Cells.Item($newRow1, 1).Value = "Languedoc.QC.1"
$ws.Cells.Item($newRow1, 2).Value = "R"
$ws.Cells.Item($newRow1, 3).Value = "Q"
$ws.Cells.Item($newRow1, 4).Value = "C"
$ws.Cells.Item($newRow1, 5).Value = "Jeu de cartes au portrait du Languedoc"
$ws.Cells.Item($newRow1, 6).Value = 1702
$ws.Cells.Item($newRow1, 7).Value = 1720
$ws.Cells.Item($newRow1, 9).Value = "Toulouse"
$ws.Cells.Item($newRow1, 10).Value = "Languedoc"
$ws.Cells.Item($newRow1, 11).Value = "Typographical letters"
$ws.Cells.Item($newRow1, 12).Value = "http://catalogue.bnf.fr/ark:/12148/cb40918053r"

# Row: Verso (2) - Queen of Clubs
$ws.Cells.Item($newRow2, 1).Value = "Languedoc.QC.2"
$ws.Cells.Item($newRow2, 2).Value = "V"
$ws.Cells.Item($newRow2, 3).Value = "Q"
$ws.Cells.Item($newRow2, 4).Value = "C"
$ws.Cells.Item($newRow2, 5).Value = "Jeu de cartes au portrait du Languedoc"
$ws.Cells.Item($newRow2, 6).Value = 1702
$ws.Cells.Item($newRow2, 7).Value = 1720
$ws.Cells.Item($newRow2, 9).Value = "Toulouse"
$ws.Cells.Item($newRow2, 10).Value = "Languedoc"
$ws.Cells.Item($newRow2, 11).Value = "Typographical letters"
$ws.Cells.Item($newRow2, 12).Value = "http://catalogue.bnf.fr/ark:/12148/cb40918053r"

# Move the active selection to reflect where editing left off.
$ws.Range("A24").Select()
